# MDSiTestResult.xlsx - "changes of 2nd May 2022"
#
# The three Job# values in column B (rows 2-4) are refreshed to the new
# job numbers produced by the latest test run:
#   B2: 32376214 -> 32378887  (H3P / FedEx)
#   B3: 32376215 -> 32378888  (H3P / UPS)
#   B4: 32376217 -> 32378889  (CPU)
#
# These values are stored as text in the original workbook (shared string
# cells with no explicit numeric formatting/style), even though they look
# like numbers. Plain "$range.Value = '...'" would let Excel auto-detect
# the numeric string and store it as a Number, so each cell is briefly
# switched to Text number-format before the assignment (forcing Excel to
# keep the literal text) and then restored to the default "Normal" style
# afterwards so no visible formatting change is left behind - matching
# the original (unstyled) B2:B4 cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "32378887"
$ws.Range("B2").Style = "Normal"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "32378888"
$ws.Range("B3").Style = "Normal"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "32378889"
$ws.Range("B4").Style = "Normal"
